$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "amhfcu_business_2023 (2)"
$wb.Worksheets.Item(3).Name = "amhfcu_personal_2023 (1)"
$wb.Worksheets.Item(4).Name = "transaction_history-2"
$wb.Worksheets.Item(5).Name = "Upwork_2023"
